$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 2059.7
$ws.Range("I34").Value = 2059.7
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2059.7
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1856.7
$ws.Range("H36").Value = 2059.7
$ws.Range("I36").Value = 2059.7
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2059.7
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1344.7
$ws.Range("H40").Value = 4399.9443
$ws.Range("I40").Value = 3280
$ws.Range("J40").Value = 5799.875
$ws.Range("K40").Value = 3280
$ws.Range("L40").Value = 5799.875
$ws.Range("M40").Value = -3105
$ws.Range("N40").Value = -6149.875
$ws.Range("H61").Value = 2398.3333
$ws.Range("I61").Value = 597.5
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 1792.5
$ws.Range("L61").Value = 18000
$ws.Range("M61").Value = -1620.5
$ws.Range("H68").Value = 110995
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 110995
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 110995
$ws.Range("N68").Value = -112493
$ws.Range("H71").Value = 110995
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 110995
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 332985
$ws.Range("N71").Value = -340473
$ws.Range("H76").Value = 3999.889
$ws.Range("I76").Value = 4399.8
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 4399.8
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -4084.8
$ws.Range("H79").Value = 3999.889
$ws.Range("I79").Value = 4399.8
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 4399.8
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -3307.8
$ws.Range("H88").Value = 584471.2
$ws.Range("I88").Value = 2450
$ws.Range("J88").Value = 681474.75
$ws.Range("K88").Value = 2450
$ws.Range("L88").Value = 681474.75
$ws.Range("M88").Value = -2044
$ws.Range("N88").Value = -682286.75
$ws.Range("H91").Value = 584471.2
$ws.Range("I91").Value = 2450
$ws.Range("J91").Value = 681474.75
$ws.Range("K91").Value = 2450
$ws.Range("L91").Value = 681474.75
$ws.Range("M91").Value = -1046
$ws.Range("N91").Value = -684282.75
$ws.Range("H133").Value = 91999.664
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 91999.664
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 91999.664
$ws.Range("N133").Value = -102119.664
$ws.Range("H137").Value = 5142.3335
$ws.Range("I137").Value = 3294.0625
$ws.Range("J137").Value = 7830.727
$ws.Range("K137").Value = 9882.1875
$ws.Range("L137").Value = 23492.181
$ws.Range("M137").Value = -7332.1875
$ws.Range("N137").Value = -28592.181
$ws.Range("H138").Value = 2979.6118
$ws.Range("I138").Value = 1830.9231
$ws.Range("J138").Value = 3187.014
$ws.Range("K138").Value = 5492.7693
$ws.Range("L138").Value = 9561.042000000001
$ws.Range("M138").Value = -352.7692999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8623915
$ws.Range("I32").Value = 9093438
$ws.Range("J32").Value = 16004.667
$ws.Range("K32").Value = 9093438
$ws.Range("L32").Value = 16004.667
$ws.Range("M32").Value = -9093151
$ws.Range("H61").Value = 9829791
$ws.Range("I61").Value = 13160994
$ws.Range("J61").Value = 92429.234
$ws.Range("K61").Value = 13160994
$ws.Range("L61").Value = 92429.234
$ws.Range("M61").Value = -13160782
$ws.Range("H94").Value = 46946
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 46946
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 46946
$ws.Range("N94").Value = -48748
$ws.Range("H97").Value = 1556.591
$ws.Range("I97").Value = 1627.35
$ws.Range("J97").Value = 849
$ws.Range("K97").Value = 1627.35
$ws.Range("L97").Value = 849
$ws.Range("M97").Value = -1131.35
$ws.Range("N97").Value = -1841
$ws.Range("H105").Value = 70500
$ws.Range("I105").Value = 30000
$ws.Range("J105").Value = 111000
$ws.Range("K105").Value = 30000
$ws.Range("L105").Value = 111000
$ws.Range("M105").Value = -26506
$ws.Range("N105").Value = -117988
$ws.Range("H136").Value = 9829791
$ws.Range("I136").Value = 13160994
$ws.Range("J136").Value = 92429.234
$ws.Range("K136").Value = 39482982
$ws.Range("L136").Value = 277287.702
$ws.Range("M136").Value = -39480432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 722.38464
$ws.Range("I22").Value = 235.72728
$ws.Range("J22").Value = 3399
$ws.Range("K22").Value = 235.72728
$ws.Range("L22").Value = 3399
$ws.Range("M22").Value = -62.72728000000001
$ws.Range("H134").Value = 717885.8
$ws.Range("I134").Value = 3098.2
$ws.Range("J134").Value = 2504854.8
$ws.Range("K134").Value = 9294.599999999999
$ws.Range("L134").Value = 7514564.399999999
$ws.Range("M134").Value = -6759.599999999999
$ws.Range("N134").Value = -7519634.399999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 973.6
$ws.Range("I16").Value = 623.3333
$ws.Range("J16").Value = 1499
$ws.Range("K16").Value = 623.3333
$ws.Range("L16").Value = 1499
$ws.Range("M16").Value = -336.3333
$ws.Range("N16").Value = -2073
$ws.Range("H31").Value = 591177.4
$ws.Range("I31").Value = 10326.2
$ws.Range("J31").Value = 1172028.5
$ws.Range("K31").Value = 10326.2
$ws.Range("L31").Value = 1172028.5
$ws.Range("M31").Value = -10031.2
$ws.Range("N31").Value = -1172618.5
$ws.Range("H34").Value = 591177.4
$ws.Range("I34").Value = 10326.2
$ws.Range("J34").Value = 1172028.5
$ws.Range("K34").Value = 10326.2
$ws.Range("L34").Value = 1172028.5
$ws.Range("M34").Value = -10124.2
$ws.Range("N34").Value = -1172432.5
$ws.Range("H58").Value = 1479.3572
$ws.Range("I58").Value = 1362.3846
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 1362.3846
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -1159.3846
$ws.Range("H110").Value = 77484.5
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 77484.5
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 77484.5
$ws.Range("N110").Value = -85664.5
$ws.Range("H113").Value = 973.6
$ws.Range("I113").Value = 623.3333
$ws.Range("J113").Value = 1499
$ws.Range("K113").Value = 623.3333
$ws.Range("L113").Value = 1499
$ws.Range("M113").Value = 1546.6667
$ws.Range("N113").Value = -5839
$ws.Range("H132").Value = 2175.0908
$ws.Range("I132").Value = 1956.6154
$ws.Range("J132").Value = 3879.2
$ws.Range("K132").Value = 5869.8462
$ws.Range("L132").Value = 11637.6
$ws.Range("M132").Value = -3339.8462
$ws.Range("N132").Value = -16697.6
$ws.Range("H136").Value = 1479.3572
$ws.Range("I136").Value = 1362.3846
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4087.1538
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1537.1538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H45").Value = 20291.334
$ws.Range("I45").Value = 842
$ws.Range("J45").Value = 30016
$ws.Range("K45").Value = 2526
$ws.Range("L45").Value = 90048
$ws.Range("M45").Value = -1994
$ws.Range("N45").Value = -91112
$ws.Range("H129").Value = 17597162
$ws.Range("I129").Value = 4160.875
$ws.Range("J129").Value = 30392070
$ws.Range("K129").Value = 12482.625
$ws.Range("L129").Value = 91176210
$ws.Range("M129").Value = -7482.625
$ws.Range("N129").Value = -91186210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 14000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 14000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 14000
$ws.Range("N18").Value = -14586
$ws.Range("H95").Value = 125040540
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 125040540
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 125040540
$ws.Range("N95").Value = -125046032
$ws.Range("H97").Value = 1727.4615
$ws.Range("I97").Value = 1845.5
$ws.Range("J97").Value = 311
$ws.Range("K97").Value = 1845.5
$ws.Range("L97").Value = 311
$ws.Range("M97").Value = -1349.5
$ws.Range("N97").Value = -1303
$ws.Range("H113").Value = 1034.1666
$ws.Range("I113").Value = 1037.2727
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1037.2727
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1132.7273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H40").Value = 3683.7856
$ws.Range("I40").Value = 2757.3
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 2757.3
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -2621.3
$ws.Range("H68").Value = 2624.5
$ws.Range("I68").Value = 2624.5
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2624.5
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1875.5
$ws.Range("H71").Value = 2624.5
$ws.Range("I71").Value = 2624.5
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 13122.5
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -9378.5
$ws.Range("H82").Value = 741.8889
$ws.Range("I82").Value = 673.44446
$ws.Range("J82").Value = 810.3333
$ws.Range("K82").Value = 673.44446
$ws.Range("L82").Value = 810.3333
$ws.Range("M82").Value = -312.44446
$ws.Range("N82").Value = -1532.3333
$ws.Range("H85").Value = 741.8889
$ws.Range("I85").Value = 673.44446
$ws.Range("J85").Value = 810.3333
$ws.Range("K85").Value = 673.44446
$ws.Range("L85").Value = 810.3333
$ws.Range("M85").Value = 574.55554
$ws.Range("N85").Value = -3306.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 45000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 45000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 45000
$ws.Range("N26").Value = -45586
$ws.Range("M26").ClearContents()
$ws.Range("H97").Value = 56411.285
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 56411.285
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 56411.285
$ws.Range("N97").Value = -58393.285
$ws.Range("H98").Value = 25882.715
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 25882.715
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 25882.715
$ws.Range("N98").Value = -31872.715
$ws.Range("H126").Value = 7424.2383
$ws.Range("I126").Value = 7433.647
$ws.Range("J126").Value = 7384.25
$ws.Range("K126").Value = 22300.941
$ws.Range("L126").Value = 22152.75
$ws.Range("M126").Value = -19830.941
$ws.Range("N126").Value = -27092.75
